$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the sheet default (unstyled) format, used to restore
# style on cells that need a forced-text (quote-prefix) assignment so no
# stray number style sticks around.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "26.893.78"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "1.730.31"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'240.55"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("D8").Value = "'0.2600"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "'0.06175"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "1.728.91"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "'16.04"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").Value = "'0.06879"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "'0.6030"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'4.463"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'77.00"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'0.9998"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "26.673.83"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "'0.9992"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "1.950.33"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'4.393"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "'8.414"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "'5.060"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "'139.72"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").Value = "'106.67"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'1.375"
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("D30").Value = "'3.949"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").Value = "'0.07926"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'3.665"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "'0.04593"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").Value = "'2.592"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").Value = "'0.6163"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "'2.462"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").Value = "'1.986"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'5.715"
$ws.Range("E41").Value = "  +5.84%  "
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "'99.93"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'0.3832"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'6.762"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").Value = "'0.1153"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").Value = "'0.05363"
$ws.Range("D48").Value = "'7.908"
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("D49").Value = "'30.12"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").Value = "'1.241"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "'51.23"
$ws.Range("E51").Value = "  -0.42%  "

# Restore default (no explicit number-format / quote-prefix) styling on the
# cells that had to be force-typed as text above.
$ws.Range("D4").Style = $defaultStyle
$ws.Range("D5").Style = $defaultStyle
$ws.Range("D6").Style = $defaultStyle
$ws.Range("D8").Style = $defaultStyle
$ws.Range("D9").Style = $defaultStyle
$ws.Range("D11").Style = $defaultStyle
$ws.Range("D12").Style = $defaultStyle
$ws.Range("D13").Style = $defaultStyle
$ws.Range("D14").Style = $defaultStyle
$ws.Range("D15").Style = $defaultStyle
$ws.Range("D16").Style = $defaultStyle
$ws.Range("D18").Style = $defaultStyle
$ws.Range("D22").Style = $defaultStyle
$ws.Range("D23").Style = $defaultStyle
$ws.Range("D24").Style = $defaultStyle
$ws.Range("D25").Style = $defaultStyle
$ws.Range("D28").Style = $defaultStyle
$ws.Range("D29").Style = $defaultStyle
$ws.Range("D30").Style = $defaultStyle
$ws.Range("D31").Style = $defaultStyle
$ws.Range("D32").Style = $defaultStyle
$ws.Range("D33").Style = $defaultStyle
$ws.Range("D34").Style = $defaultStyle
$ws.Range("D35").Style = $defaultStyle
$ws.Range("D36").Style = $defaultStyle
$ws.Range("D38").Style = $defaultStyle
$ws.Range("D39").Style = $defaultStyle
$ws.Range("D41").Style = $defaultStyle
$ws.Range("D43").Style = $defaultStyle
$ws.Range("D44").Style = $defaultStyle
$ws.Range("D45").Style = $defaultStyle
$ws.Range("D46").Style = $defaultStyle
$ws.Range("D47").Style = $defaultStyle
$ws.Range("D48").Style = $defaultStyle
$ws.Range("D49").Style = $defaultStyle
$ws.Range("D50").Style = $defaultStyle
$ws.Range("D51").Style = $defaultStyle
